# Insert a new weekly price record for "Vega Modelo de Temuco" (Jengibre)
# above the current row 288 — this shifts the existing rows 288-331 down
# to 289-332 and grows the used range to A1:R332.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(288).Insert()

$ws.Cells.Item(288, 1).Value  = 10
$ws.Cells.Item(288, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(288, 3).Value  = "La Araucanía"
$ws.Cells.Item(288, 4).Value  = 45131
$ws.Cells.Item(288, 5).Value  = 9
$ws.Cells.Item(288, 6).Value  = 100114007
$ws.Cells.Item(288, 7).Value  = "Jengibre"
$ws.Cells.Item(288, 8).Value  = "Sin especificar"
$ws.Cells.Item(288, 9).Value  = "Primera"
$ws.Cells.Item(288, 10).Value = 35
$ws.Cells.Item(288, 11).Value = 20000
$ws.Cells.Item(288, 12).Value = 20000
$ws.Cells.Item(288, 13).Value = 20000
$ws.Cells.Item(288, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(288, 15).Value = "Perú"
$ws.Cells.Item(288, 16).Value = 1538
$ws.Cells.Item(288, 17).Value = 13
$ws.Cells.Item(288, 18).Value = "Hortaliza"
